$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Rushing")
$ws2 = $wb.Worksheets.Item("Receiving")

# --- Sheet1 (Rushing): Week 15 stat updates ---
# Row 4: E.Elliott
$ws1.Cells.Item(4, 3).Value = 120
$ws1.Cells.Item(4, 4).Value = 56
$ws1.Cells.Item(4, 5).Value = 25
$ws1.Cells.Item(4, 6).Value = 35

# Row 5: T.Pollard
$ws1.Cells.Item(5, 3).Value = 72
$ws1.Cells.Item(5, 4).Value = 39

# --- Sheet2 (Receiving): Week 15 stat updates ---
# Row 2: E.Elliott
$ws2.Cells.Item(2, 3).Value = 44
$ws2.Cells.Item(2, 4).Value = 32
$ws2.Cells.Item(2, 7).Value = 16
$ws2.Cells.Item(2, 8).Value = 11

# Row 3: T.Pollard
$ws2.Cells.Item(3, 3).Value = 29
$ws2.Cells.Item(3, 4).Value = 25

# Row 5: A.Cooper
$ws2.Cells.Item(5, 3).Value = 56
$ws2.Cells.Item(5, 4).Value = 41
$ws2.Cells.Item(5, 5).Value = 22
$ws2.Cells.Item(5, 7).Value = 13
$ws2.Cells.Item(5, 8).Value = 9

# Row 6: C.Lamb
$ws2.Cells.Item(6, 3).Value = 73
$ws2.Cells.Item(6, 4).Value = 51
$ws2.Cells.Item(6, 5).Value = 32
$ws2.Cells.Item(6, 7).Value = 11

# Row 7: M.Gallup
$ws2.Cells.Item(7, 3).Value = 41
$ws2.Cells.Item(7, 4).Value = 33
$ws2.Cells.Item(7, 5).Value = 12

# Row 10: M.Turner
$ws2.Cells.Item(10, 3).Value = 9
$ws2.Cells.Item(10, 4).Value = 8

# Row 12: D.Schultz
$ws2.Cells.Item(12, 3).Value = 59
$ws2.Cells.Item(12, 4).Value = 47
$ws2.Cells.Item(12, 7).Value = 9
$ws2.Cells.Item(12, 8).Value = 5

# --- Sheet2: new Week 16 simulated player row (row 14) ---
$ws2.Cells.Item(14, 1).Value = 12
$ws2.Cells.Item(14, 2).Value = "J.Sprinkle"
$ws2.Cells.Item(14, 3).Value = 2
$ws2.Cells.Item(14, 4).Value = 2
$ws2.Cells.Item(14, 5).Value = 0
$ws2.Cells.Item(14, 6).Value = 0
$ws2.Cells.Item(14, 7).Value = 1
$ws2.Cells.Item(14, 8).Value = 1

# Copy the row-label style (bold, centered, bordered) from A13 onto the new A14 cell
$ws2.Range("A13").Copy()
$ws2.Range("A14").PasteSpecial(-4122)
